# Fix: "flat error" — restructure the workbook:
#   1. Insert a brand-new first sheet "시설물 현황" (facility overview) with a small
#      formatted summary table.
#   2. Keep the old empty "Sheet" (now 2nd tab).
#   3. Trim the "손상현황표" (damage-status table) down to a single data row and
#      update its print area accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "시설물 현황" sheet and move it to the very first position.
# ---------------------------------------------------------------------------
$infoSheet = $wb.Worksheets.Add()
$infoSheet.Name = "시설물 현황"
$infoSheet.Move($wb.Worksheets.Item(1))

# Match the page margins used throughout the rest of the workbook (values are
# expressed in points: 0.75in = 54pt, 1in = 72pt, 0.5in = 36pt).
$infoSheet.PageSetup.LeftMargin = 54
$infoSheet.PageSetup.RightMargin = 54
$infoSheet.PageSetup.TopMargin = 72
$infoSheet.PageSetup.BottomMargin = 72
$infoSheet.PageSetup.HeaderMargin = 36
$infoSheet.PageSetup.FooterMargin = 36

# Base style for the whole used range: centered both ways.
$full = $infoSheet.Range("A1:G10")
$full.HorizontalAlignment = -4108
$full.VerticalAlignment = -4108

# Label cells get a light-gray fill on top of the centered alignment.
$labelRanges = "B4", "E4", "B5", "E5", "B6", "E6", "B7", "E7", "B8", "D8", "F8", "B9"
foreach ($addr in $labelRanges) {
    $infoSheet.Range($addr).Interior.Color = 13421772
}

# 준공일자 (completion date) value gets a date number format.
$infoSheet.Range("F5").NumberFormat = "yyyy-mm-dd"

# --- Section headings -------------------------------------------------------
$infoSheet.Range("B2").Value = "□ 시설물 현황"
$infoSheet.Range("B3").Value = "가. 일반현황"

# --- Row 4 -------------------------------------------------------------------
$infoSheet.Range("B4").Value = "시설물명"
$infoSheet.Range("C4").Value = "정부춘천"
$infoSheet.Range("E4").Value = "시설물번호"
$infoSheet.Range("F4").Value = "AR2003-0009512"

# --- Row 5 -------------------------------------------------------------------
$infoSheet.Range("B5").Value = "시설물위치"
$infoSheet.Range("C5").Value = "AR2003-0009512"
$infoSheet.Range("E5").Value = "준공일자"
$infoSheet.Range("F5").Value = 44559

# --- Row 6 -------------------------------------------------------------------
$infoSheet.Range("B6").Value = "용도"
$infoSheet.Range("C6").Value = "공공업무"
$infoSheet.Range("E6").Value = "시설물규모"
$infoSheet.Range("F6").Value = "지하1층"

# --- Row 7 -------------------------------------------------------------------
$infoSheet.Range("B7").Value = "구조형식"
$infoSheet.Range("C7").Value = "철근"
$infoSheet.Range("E7").Value = "부대시설"

# --- Row 8 -------------------------------------------------------------------
$infoSheet.Range("B8").Value = "종별"
$infoSheet.Range("C8").Value = "3층"
$infoSheet.Range("D8").Value = "전차안전등급"
$infoSheet.Range("E8").Value = "b등급"
$infoSheet.Range("F8").Value = "점검결과안전등급"
$infoSheet.Range("G8").Value = "b등급(8.7점)"

# --- Row 9 / 10 ----------------------------------------------------------------
$infoSheet.Range("B9").Value = "규모 및 제원 추가사항"
$infoSheet.Range("B10").Value = "없음"

# --- Merges --------------------------------------------------------------------
$infoSheet.Range("C4:D4").Merge()
$infoSheet.Range("C5:D5").Merge()
$infoSheet.Range("C6:D6").Merge()
$infoSheet.Range("C7:D7").Merge()
$infoSheet.Range("F4:G4").Merge()
$infoSheet.Range("F5:G5").Merge()
$infoSheet.Range("F6:G6").Merge()
$infoSheet.Range("F7:G7").Merge()
$infoSheet.Range("B9:G9").Merge()
$infoSheet.Range("B10:G10").Merge()

# ---------------------------------------------------------------------------
# 2. Trim "손상현황표" down to a single damage record and refresh its data.
# ---------------------------------------------------------------------------
$damageSheet = $wb.Worksheets.Item("손상현황표")

# Drop row 4 (the second damage record) — row 3 becomes the only data row.
$damageSheet.Rows(4).Delete()

# Update row 3 to reflect the single remaining damage record.
$damageSheet.Range("B3").Value = "5층"
$damageSheet.Range("C3").Value = "옥탑층"
$damageSheet.Range("D3").Value = "벽재"
$damageSheet.Range("E3").Value = "0.2x3500 수직, 수평"
$damageSheet.Range("F3").Value = "'3"
$damageSheet.Range("G3").Value = "'2"
$damageSheet.Range("H3").Value = "신규"
$damageSheet.Range("I3").Value = "모름"
$damageSheet.Range("J3").Value = ""

# Print area now covers only the remaining two rows of the table.
$damageSheet.PageSetup.PrintArea = '$A$2:$J$3'

Write-Output "done"
